# Apply the updated "人气/热度" (F column) values produced by the latest
# data-refresh run ("Update gh-pages to output generated at 456a3b4").
# Sheet "展览" and the combined "全部类型" sheet each contain a row per
# exhibition event; both copies of the F-column counter need to be bumped
# to match the newly generated numbers.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAll        = $wb.Worksheets.Item("全部类型")

# Cell -> new value updates for the "展览" sheet
$exhibitionUpdates = @{
    "F2"  = 284
    "F4"  = 292
    "F9"  = 1561
    "F13" = 2615
    "F16" = 6837
    "F18" = 6982
    "F20" = 2257
    "F21" = 3042
    "F23" = 211
    "F24" = 110
    "F25" = 1754
    "F27" = 285
    "F28" = 860
    "F30" = 164
    "F31" = 28
    "F32" = 364
    "F34" = 2393
    "F38" = 935
    "F39" = 194
    "F40" = 437
}

foreach ($cell in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range($cell).Value = $exhibitionUpdates[$cell]
}

# Cell -> new value updates for the "全部类型" sheet
$allTypesUpdates = @{
    "F4"  = 284
    "F6"  = 292
    "F10" = 1561
    "F15" = 2615
    "F21" = 6837
    "F23" = 6982
    "F25" = 2257
    "F26" = 3042
    "F29" = 211
    "F32" = 1754
    "F35" = 285
    "F36" = 860
    "F37" = 164
    "F38" = 28
    "F39" = 364
    "F41" = 2393
    "F46" = 935
    "F47" = 194
    "F48" = 437
}

foreach ($cell in $allTypesUpdates.Keys) {
    $sheetAll.Range($cell).Value = $allTypesUpdates[$cell]
}
